# Loan RBI, Variable Instalments
# Insert a new (empty) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / heading / "Outstanding" columns one
# place to the right, then leave the selection on the newly widened sheet
# (which makes it the active sheet/tab instead of "Transactions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a whole new column at N - this shifts N:P -> O:Q automatically,
# copying cell styles/values along with them.
$ws.Columns("N:N").Insert()

# Give the freshly inserted column a plain custom width of 11 (no bestFit),
# matching the author's manual resize after inserting the column.
$ws.Columns("N:N").ColumnWidth = 10.17

# Select R3 on the Repayment schedule sheet - this both matches the final
# selection recorded in the workbook and makes this sheet the active tab
# (clearing tabSelected on the previously-active "Transactions" sheet).
$ws.Range("R3").Select() | Out-Null
